$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 974, shifting existing rows 974-1000 down to 975-1001
$ws.Rows.Item(974).Insert()

# Fill in the new row 974 with the new contest data
$ws.Cells.Item(974, 1).Value = 3588
$ws.Cells.Item(974, 2).Value = "15/01/2026"
$ws.Cells.Item(974, 3).Value = 3
$ws.Cells.Item(974, 4).Value = 5
$ws.Cells.Item(974, 5).Value = 7
$ws.Cells.Item(974, 6).Value = 8
$ws.Cells.Item(974, 7).Value = 9
$ws.Cells.Item(974, 8).Value = 11
$ws.Cells.Item(974, 9).Value = 14
$ws.Cells.Item(974, 10).Value = 15
$ws.Cells.Item(974, 11).Value = 16
$ws.Cells.Item(974, 12).Value = 17
$ws.Cells.Item(974, 13).Value = 19
$ws.Cells.Item(974, 14).Value = 21
$ws.Cells.Item(974, 15).Value = 22
$ws.Cells.Item(974, 16).Value = 23
$ws.Cells.Item(974, 17).Value = 24
